$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.734.07'
$ws.Range("E2").Value = '  -1.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.760.42'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.90'
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4434'
$ws.Range("E7").Value = '  -2.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3740'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.44'
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07758'
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.129'
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.79'
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.204'
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.377'
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.760.25'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.35'
$ws.Range("E17").Value = '  +12.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001082'
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06243'
$ws.Range("E19").Value = '  -7.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.44'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.197'
$ws.Range("E22").Value = '  -2.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5329'
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.767.38'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.65'
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.326'
$ws.Range("E26").Value = '  -4.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.78'
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.76'
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.369'
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.960.70'
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.33'
$ws.Range("E31").Value = '  -2.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.218'
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.793'
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09293'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.656'
$ws.Range("E35").Value = '  -9.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.72'
$ws.Range("E36").Value = '  +4.90%  '
$ws.Range("E37").Value = '  -7.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02336'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06159'
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6514'
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.105'
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.198'
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.034'
$ws.Range("E43").Value = '  -3.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.416'
$ws.Range("E44").Value = '  -4.33%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.84'
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6036'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.766'
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.35'
$ws.Range("E49").Value = '  -2.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.002'
$ws.Range("E50").Value = '  -1.51%  '
$ws.Range("E51").Value = '  -1.23%  '
